$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# The item "سرنجات 3 سم" (row 18, serial #15) is removed from the list.
# Deleting the entire row shifts every row below it up by one (rows 19-22 -> 18-21),
# including the merged cell ranges.
$ws.Rows("18").Delete()

# After the shift, the serial-number column (A) and the row heights for the two
# rows that slid up need to go back to their own original values (they are not
# part of the shift in the source workbook - only the item details moved up).
$ws.Range("A18").Value = 15
$ws.Range("A19").Value = 16
$ws.Rows("18").RowHeight = 25.5
$ws.Rows("19").RowHeight = 24.75

# Row 20 now holds the totals row (previously row 21). Its total must reflect
# the removal of the deleted item's value (-10), i.e. 209.6 - (-10) = 219.6,
# and its row height grows slightly to match the new layout.
$ws.Range("K20").Value = 219.6
$ws.Rows("20").RowHeight = 26.25

Write-Output "Row for 'سرنجات 3 سم' removed; totals and layout updated."
